$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the extra duplicate columns U:AD (data no longer needed)
$ws.Range("U1:AD19").EntireColumn.Delete()

# Step 2: dereference the shared strings that need reordering/insertion so the
# rebuilt sharedStrings table lands in the target order once we re-set them below.
$ws.Range("B16:B19").ClearContents()
$ws.Range("C2:J2").ClearContents()
$ws.Range("K2:T2").ClearContents()

# Step 3a: Holden rows (reuses rows 16-19, previously HexGrid) -- new shared strings,
# inserted first so they land right after OffsetATD in the table.
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("B17").Value = "Holden5"
$ws.Range("B18").Value = "Holden10"
$ws.Range("B19").Value = "Holden15"

# Step 3b: re-add the HexGrid rows as new rows 20-23, copying column-A number format
# from row 19 so the bold/centered/bordered style (s=1) carries over like Excel would
# for a freshly extended table.
$ws.Range("A19").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21
$ws.Range("B20").Value = "HexGrid-90degTilt2.5degRes"
$ws.Range("B21").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B22").Value = "HexGrid-90degTilt10degRes"
$ws.Range("B23").Value = "HexGrid-90degTilt15degRes"

# Step 3c: HKL column header labels (row 2), reordered
$ws.Range("C2").Value = "[3, 2, 1]"
$ws.Range("D2").Value = "[2, 2, 2]"
$ws.Range("E2").Value = "[3, 1, 0]"
$ws.Range("F2").Value = "[1, 1, 0]"
$ws.Range("G2").Value = "[2, 2, 0]"
$ws.Range("H2").Value = "[2, 0, 0]"
$ws.Range("I2").Value = "[4, 0, 0]"
$ws.Range("J2").Value = "[2, 1, 1]"

# Step 3d: pair-count column header labels (row 2), re-added last
$ws.Range("K2").Value = "1Pair-A"
$ws.Range("L2").Value = "1Pair-B"
$ws.Range("M2").Value = "2Pairs-A"
$ws.Range("N2").Value = "2Pairs-B"
$ws.Range("O2").Value = "3Pairs-A"
$ws.Range("P2").Value = "3Pairs-B"
$ws.Range("Q2").Value = "3Pairs-C"
$ws.Range("R2").Value = "4Pairs"
$ws.Range("S2").Value = "5A4F"
$ws.Range("T2").Value = "MaxUnique"

# Step 4: numeric data rows 3-19 (existing rows) -- HKL columns permuted + refreshed values
# Row 3: BT8Hex_2.5
$ws.Range("C3").Value = 0.9992840340062283
$ws.Range("D3").Value = 0.9878296850411552
$ws.Range("E3").Value = 1.001433247393378
$ws.Range("F3").Value = 1.010313643421926
$ws.Range("G3").Value = 1.010313643421926
$ws.Range("H3").Value = 0.9988668487488201
$ws.Range("I3").Value = 0.9988668487488201
$ws.Range("J3").Value = 0.9957608877534186
$ws.Range("K3").Value = 1.010313643421926
$ws.Range("L3").Value = 0.9957608877534186
$ws.Range("M3").Value = 0.9973138682511193
$ws.Range("N3").Value = 0.9973138682511193
$ws.Range("O3").Value = 0.9986869946318722
$ws.Range("P3").Value = 1.001647126641388
$ws.Range("Q3").Value = 1.001647126641388
$ws.Range("R3").Value = 1.003813755836523
$ws.Range("S3").Value = 1.003813755836523
$ws.Range("T3").Value = 0.9989147243941544

# Row 4: BT8Hex_5
$ws.Range("C4").Value = 0.9986395087887119
$ws.Range("D4").Value = 0.9761204295953141
$ws.Range("E4").Value = 1.002724914057606
$ws.Range("F4").Value = 1.020327310461459
$ws.Range("G4").Value = 1.020327310461459
$ws.Range("H4").Value = 0.997589033456908
$ws.Range("I4").Value = 0.997589033456908
$ws.Range("J4").Value = 0.9916335107254579
$ws.Range("K4").Value = 1.020327310461459
$ws.Range("L4").Value = 0.9916335107254579
$ws.Range("M4").Value = 0.9946112720911829
$ws.Range("N4").Value = 0.9946112720911829
$ws.Range("O4").Value = 0.9973158194133239
$ws.Range("P4").Value = 1.003183284881275
$ws.Range("Q4").Value = 1.003183284881275
$ws.Range("R4").Value = 1.007469291276321
$ws.Range("S4").Value = 1.007469291276321
$ws.Range("T4").Value = 0.9978391178475761

# Row 5: BT8Hex_10
$ws.Range("C5").Value = 0.997194255252568
$ws.Range("D5").Value = 0.9570653835116948
$ws.Range("E5").Value = 1.005238073061536
$ws.Range("F5").Value = 1.036544086658786
$ws.Range("G5").Value = 1.036544086658786
$ws.Range("H5").Value = 0.9971226400129226
$ws.Range("I5").Value = 0.9971226400129226
$ws.Range("J5").Value = 0.9850024173647374
$ws.Range("K5").Value = 1.036544086658786
$ws.Range("L5").Value = 0.9850024173647374
$ws.Range("M5").Value = 0.9910625286888299
$ws.Range("N5").Value = 0.9910625286888299
$ws.Range("O5").Value = 0.9957877101463986
$ws.Range("P5").Value = 1.006223048012149
$ws.Range("Q5").Value = 1.006223048012149
$ws.Range("R5").Value = 1.013803307673808
$ws.Range("S5").Value = 1.013803307673808
$ws.Range("T5").Value = 0.9963611426437073

# Row 6: BT8Hex_15
$ws.Range("C6").Value = 0.9950493389529582
$ws.Range("D6").Value = 0.9395788043893717
$ws.Range("E6").Value = 1.007758021776985
$ws.Range("F6").Value = 1.054306147978222
$ws.Range("G6").Value = 1.054306147978222
$ws.Range("H6").Value = 0.9969158964927514
$ws.Range("I6").Value = 0.9969158964927514
$ws.Range("J6").Value = 0.9799850216717438
$ws.Range("K6").Value = 1.054306147978222
$ws.Range("L6").Value = 0.9799850216717438
$ws.Range("M6").Value = 0.9884504590822476
$ws.Range("N6").Value = 0.9884504590822476
$ws.Range("O6").Value = 0.9948863133138269
$ws.Range("P6").Value = 1.010402355380906
$ws.Range("Q6").Value = 1.010402355380906
$ws.Range("R6").Value = 1.021378303530235
$ws.Range("S6").Value = 1.021378303530235
$ws.Range("T6").Value = 0.9955988718770055

# Row 7: Spiral2.5
$ws.Range("C7").Value = 0.9998830707076818
$ws.Range("D7").Value = 1.000694579872367
$ws.Range("E7").Value = 1.000019750360145
$ws.Range("F7").Value = 0.9995486872504216
$ws.Range("G7").Value = 0.9995486872504216
$ws.Range("H7").Value = 1.000784820436869
$ws.Range("I7").Value = 1.000784820436869
$ws.Range("J7").Value = 1.000124776338551
$ws.Range("K7").Value = 0.9995486872504216
$ws.Range("L7").Value = 1.000124776338551
$ws.Range("M7").Value = 1.00045479838771
$ws.Range("N7").Value = 1.00045479838771
$ws.Range("O7").Value = 1.000309782378522
$ws.Range("P7").Value = 1.000152761341947
$ws.Range("Q7").Value = 1.000152761341947
$ws.Range("R7").Value = 1.000001742819066
$ws.Range("S7").Value = 1.000001742819066
$ws.Range("T7").Value = 1.000175947494339

# Row 8: Spiral5
$ws.Range("C8").Value = 0.999791856583333
$ws.Range("D8").Value = 1.001586014349713
$ws.Range("E8").Value = 0.9998161643941119
$ws.Range("F8").Value = 0.9996113584745954
$ws.Range("G8").Value = 0.9996113584745954
$ws.Range("H8").Value = 1.001386073157659
$ws.Range("I8").Value = 1.001386073157659
$ws.Range("J8").Value = 1.000170066611662
$ws.Range("K8").Value = 0.9996113584745954
$ws.Range("L8").Value = 1.000170066611662
$ws.Range("M8").Value = 1.000778069884661
$ws.Range("N8").Value = 1.000778069884661
$ws.Range("O8").Value = 1.000457434721144
$ws.Range("P8").Value = 1.000389166081306
$ws.Range("Q8").Value = 1.000389166081306
$ws.Range("R8").Value = 1.000194714179628
$ws.Range("S8").Value = 1.000194714179628
$ws.Range("T8").Value = 1.000393588928512

# Row 9: Spiral7.5
$ws.Range("C9").Value = 0.999776081333195
$ws.Range("D9").Value = 1.001184040744556
$ws.Range("E9").Value = 0.9996084037569631
$ws.Range("F9").Value = 1.00078217216279
$ws.Range("G9").Value = 1.00078217216279
$ws.Range("H9").Value = 1.001351079669106
$ws.Range("I9").Value = 1.001351079669106
$ws.Range("J9").Value = 0.9998220043066715
$ws.Range("K9").Value = 1.00078217216279
$ws.Range("L9").Value = 0.9998220043066715
$ws.Range("M9").Value = 1.000586541987889
$ws.Range("N9").Value = 1.000586541987889
$ws.Range("O9").Value = 1.000260495910914
$ws.Range("P9").Value = 1.00065175204619
$ws.Range("Q9").Value = 1.000651752046189
$ws.Range("R9").Value = 1.00068435707534
$ws.Range("S9").Value = 1.00068435707534
$ws.Range("T9").Value = 1.00042063032888

# Row 10: Spiral10
$ws.Range("C10").Value = 0.9994253486345679
$ws.Range("D10").Value = 1.003902981967662
$ws.Range("E10").Value = 0.9993130034894974
$ws.Range("F10").Value = 1.00004921021417
$ws.Range("G10").Value = 1.00004921021417
$ws.Range("H10").Value = 1.003626066628297
$ws.Range("I10").Value = 1.003626066628297
$ws.Range("J10").Value = 1.000137055221159
$ws.Range("K10").Value = 1.00004921021417
$ws.Range("L10").Value = 1.000137055221159
$ws.Range("M10").Value = 1.001881560924728
$ws.Range("N10").Value = 1.001881560924728
$ws.Range("O10").Value = 1.001025375112985
$ws.Range("P10").Value = 1.001270777354542
$ws.Range("Q10").Value = 1.001270777354542
$ws.Range("R10").Value = 1.000965385569449
$ws.Range("S10").Value = 1.000965385569449
$ws.Range("T10").Value = 1.001075611025892

# Row 11: Spiral15
$ws.Range("C11").Value = 0.9991909099803656
$ws.Range("D11").Value = 1.003698294341228
$ws.Range("E11").Value = 0.9984249087532314
$ws.Range("F11").Value = 1.004327881850952
$ws.Range("G11").Value = 1.004327881850952
$ws.Range("H11").Value = 1.004400821146334
$ws.Range("I11").Value = 1.004400821146334
$ws.Range("J11").Value = 0.998810231350628
$ws.Range("K11").Value = 1.004327881850952
$ws.Range("L11").Value = 0.998810231350628
$ws.Range("M11").Value = 1.001605526248481
$ws.Range("N11").Value = 1.001605526248481
$ws.Range("O11").Value = 1.000545320416731
$ws.Range("P11").Value = 1.002512978115971
$ws.Range("Q11").Value = 1.002512978115971
$ws.Range("R11").Value = 1.002966704049716
$ws.Range("S11").Value = 1.002966704049716
$ws.Range("T11").Value = 1.00147550790379

# Row 12: OffsetF45
$ws.Range("C12").Value = 1.010562656342278
$ws.Range("D12").Value = 0.6076480114475263
$ws.Range("E12").Value = 1.042499214400601
$ws.Range("F12").Value = 1.226152087430778
$ws.Range("G12").Value = 1.226152087430778
$ws.Range("H12").Value = 0.8135648193026975
$ws.Range("I12").Value = 0.8135648193026975
$ws.Range("J12").Value = 0.9195313905494017
$ws.Range("K12").Value = 1.226152087430778
$ws.Range("L12").Value = 0.9195313905494017
$ws.Range("M12").Value = 0.8665481049260496
$ws.Range("N12").Value = 0.8665481049260496
$ws.Range("O12").Value = 0.9251984747509002
$ws.Range("P12").Value = 0.9864160990942925
$ws.Range("Q12").Value = 0.9864160990942925
$ws.Range("R12").Value = 1.046350096178414
$ws.Range("S12").Value = 1.046350096178414
$ws.Range("T12").Value = 0.9366596965788805

# Row 13: OffsetA45
$ws.Range("C13").Value = 0.9769692938969803
$ws.Range("D13").Value = 1.146773738040237
$ws.Range("E13").Value = 1.001305353144287
$ws.Range("F13").Value = 0.9259020730545918
$ws.Range("G13").Value = 0.9259020730545918
$ws.Range("H13").Value = 1.19622858816455
$ws.Range("I13").Value = 1.19622858816455
$ws.Range("J13").Value = 0.9980379874609863
$ws.Range("K13").Value = 0.9259020730545918
$ws.Range("L13").Value = 0.9980379874609863
$ws.Range("M13").Value = 1.097133287812768
$ws.Range("N13").Value = 1.097133287812768
$ws.Range("O13").Value = 1.065190642923275
$ws.Range("P13").Value = 1.04005621622671
$ws.Range("Q13").Value = 1.04005621622671
$ws.Range("R13").Value = 1.01151768043368
$ws.Range("S13").Value = 1.01151768043368
$ws.Range("T13").Value = 1.040869505626939

# Row 14: OffsetFTD
$ws.Range("C14").Value = 0.9279295441835341
$ws.Range("D14").Value = 1.166244549060898
$ws.Range("E14").Value = 1.101895647005625
$ws.Range("F14").Value = 0.7594890053520396
$ws.Range("G14").Value = 0.7594890053520396
$ws.Range("H14").Value = 1.440559545107162
$ws.Range("I14").Value = 1.440559545107162
$ws.Range("J14").Value = 0.992844387141834
$ws.Range("K14").Value = 0.7594890053520396
$ws.Range("L14").Value = 0.992844387141834
$ws.Range("M14").Value = 1.216701966124498
$ws.Range("N14").Value = 1.216701966124498
$ws.Range("O14").Value = 1.178433193084874
$ws.Range("P14").Value = 1.064297645867012
$ws.Range("Q14").Value = 1.064297645867012
$ws.Range("R14").Value = 0.9880954857382689
$ws.Range("S14").Value = 0.9880954857382689
$ws.Range("T14").Value = 1.064827112975182

# Row 15: OffsetATD
$ws.Range("C15").Value = 1.003399005809783
$ws.Range("D15").Value = 0.8368033757414323
$ws.Range("E15").Value = 1.040837308069024
$ws.Range("F15").Value = 1.018211072495288
$ws.Range("G15").Value = 1.018211072495288
$ws.Range("H15").Value = 0.8597309207966093
$ws.Range("I15").Value = 0.8597309207966093
$ws.Range("J15").Value = 0.9981067985513478
$ws.Range("K15").Value = 1.018211072495288
$ws.Range("L15").Value = 0.9981067985513478
$ws.Range("M15").Value = 0.9289188596739786
$ws.Range("N15").Value = 0.9289188596739786
$ws.Range("O15").Value = 0.9662250091389938
$ws.Range("P15").Value = 0.958682930614415
$ws.Range("Q15").Value = 0.9586829306144149
$ws.Range("R15").Value = 0.973564966084633
$ws.Range("S15").Value = 0.973564966084633
$ws.Range("T15").Value = 0.9595147469105806

# Row 16: Holden2.5
$ws.Range("C16").Value = 0.9861619971661297
$ws.Range("D16").Value = 0.7860611369087336
$ws.Range("E16").Value = 1.029578190031759
$ws.Range("F16").Value = 1.171200274466133
$ws.Range("G16").Value = 1.171200274466133
$ws.Range("H16").Value = 0.9919431473805873
$ws.Range("I16").Value = 0.9919431473805873
$ws.Range("J16").Value = 0.9252190678924933
$ws.Range("K16").Value = 1.171200274466133
$ws.Range("L16").Value = 0.9252190678924933
$ws.Range("M16").Value = 0.9585811076365403
$ws.Range("N16").Value = 0.9585811076365403
$ws.Range("O16").Value = 0.9822468017682798
$ws.Range("P16").Value = 1.029454163246404
$ws.Range("Q16").Value = 1.029454163246404
$ws.Range("R16").Value = 1.064890691051336
$ws.Range("S16").Value = 1.064890691051336
$ws.Range("T16").Value = 0.981693968974306

# Row 17: Holden5
$ws.Range("C17").Value = 0.9859476713722346
$ws.Range("D17").Value = 0.8795238349883636
$ws.Range("E17").Value = 1.025365907998686
$ws.Range("F17").Value = 1.088111623563394
$ws.Range("G17").Value = 1.088111623563394
$ws.Range("H17").Value = 1.024812064209199
$ws.Range("I17").Value = 1.024812064209199
$ws.Range("J17").Value = 0.9553944319387363
$ws.Range("K17").Value = 1.088111623563394
$ws.Range("L17").Value = 0.9553944319387363
$ws.Range("M17").Value = 0.9901032480739678
$ws.Range("N17").Value = 0.9901032480739678
$ws.Range("O17").Value = 1.001857468048874
$ws.Range("P17").Value = 1.022772706570443
$ws.Range("Q17").Value = 1.022772706570443
$ws.Range("R17").Value = 1.039107435818681
$ws.Range("S17").Value = 1.039107435818681
$ws.Range("T17").Value = 0.9931925890117689

# Row 18: Holden10
$ws.Range("C18").Value = 0.9857573809560516
$ws.Range("D18").Value = 1.064379621914434
$ws.Range("E18").Value = 1.016501908860489
$ws.Range("F18").Value = 0.9248606186270631
$ws.Range("G18").Value = 0.9248606186270631
$ws.Range("H18").Value = 1.088668912646401
$ws.Range("I18").Value = 1.088668912646401
$ws.Range("J18").Value = 1.014949669110977
$ws.Range("K18").Value = 0.9248606186270631
$ws.Range("L18").Value = 1.014949669110977
$ws.Range("M18").Value = 1.051809290878689
$ws.Range("N18").Value = 1.051809290878689
$ws.Range("O18").Value = 1.040040163539289
$ws.Range("P18").Value = 1.009493066794813
$ws.Range("Q18").Value = 1.009493066794813
$ws.Range("R18").Value = 0.9883349547528759
$ws.Range("S18").Value = 0.9883349547528759
$ws.Range("T18").Value = 1.015853018685902

# Row 19: Holden15
$ws.Range("C19").Value = 0.982566848064498
$ws.Range("D19").Value = 1.104062741632205
$ws.Range("E19").Value = 1.016959225228426
$ws.Range("F19").Value = 0.8916203183718099
$ws.Range("G19").Value = 0.8916203183718099
$ws.Range("H19").Value = 1.116877088677985
$ws.Range("I19").Value = 1.116877088677985
$ws.Range("J19").Value = 1.025466372349891
$ws.Range("K19").Value = 0.8916203183718099
$ws.Range("L19").Value = 1.025466372349891
$ws.Range("M19").Value = 1.071171730513938
$ws.Range("N19").Value = 1.071171730513938
$ws.Range("O19").Value = 1.053100895418767
$ws.Range("P19").Value = 1.011321259799895
$ws.Range("Q19").Value = 1.011321259799895
$ws.Range("R19").Value = 0.981396024442874
$ws.Range("S19").Value = 0.981396024442874
$ws.Range("T19").Value = 1.022925432387469

# Step 5: numeric data rows 20-23 (new rows, former HexGrid data, HKL columns permuted)
# Row 20: HexGrid-90degTilt2.5degRes
$ws.Range("C20").Value = 0.9999309407965052
$ws.Range("D20").Value = 1.00304405912742
$ws.Range("E20").Value = 1.000135616856035
$ws.Range("F20").Value = 0.996326224979731
$ws.Range("G20").Value = 0.996326224979731
$ws.Range("H20").Value = 1.000918552906954
$ws.Range("I20").Value = 1.000918552906954
$ws.Range("J20").Value = 1.001473472435609
$ws.Range("K20").Value = 0.996326224979731
$ws.Range("L20").Value = 1.001473472435609
$ws.Range("M20").Value = 1.001196012671282
$ws.Range("N20").Value = 1.001196012671282
$ws.Range("O20").Value = 1.000842547399533
$ws.Range("P20").Value = 0.9995727501074314
$ws.Range("Q20").Value = 0.9995727501074313
$ws.Range("R20").Value = 0.9987611188255062
$ws.Range("S20").Value = 0.9987611188255062
$ws.Range("T20").Value = 1.000304811183709

# Row 21: HexGrid-90degTilt5degRes
$ws.Range("C21").Value = 0.9998899103356096
$ws.Range("D21").Value = 0.9968942296842329
$ws.Range("E21").Value = 0.9999635757726629
$ws.Range("F21").Value = 1.004022653101704
$ws.Range("G21").Value = 1.004022653101704
$ws.Range("H21").Value = 0.9997593175139168
$ws.Range("I21").Value = 0.9997593175139168
$ws.Range("J21").Value = 0.9982742615212041
$ws.Range("K21").Value = 1.004022653101704
$ws.Range("L21").Value = 0.9982742615212041
$ws.Range("M21").Value = 0.9990167895175605
$ws.Range("N21").Value = 0.9990167895175605
$ws.Range("O21").Value = 0.9993323849359279
$ws.Range("P21").Value = 1.000685410712275
$ws.Range("Q21").Value = 1.000685410712275
$ws.Range("R21").Value = 1.001519721309632
$ws.Range("S21").Value = 1.001519721309632
$ws.Range("T21").Value = 0.9998006579882217

# Row 22: HexGrid-90degTilt10degRes
$ws.Range("C22").Value = 0.999384395862361
$ws.Range("D22").Value = 1.004031061709294
$ws.Range("E22").Value = 1.000651426531763
$ws.Range("F22").Value = 0.9948425901970857
$ws.Range("G22").Value = 0.9948425901970857
$ws.Range("H22").Value = 1.000731750527894
$ws.Range("I22").Value = 1.000731750527894
$ws.Range("J22").Value = 1.002561852172989
$ws.Range("K22").Value = 0.9948425901970857
$ws.Range("L22").Value = 1.002561852172989
$ws.Range("M22").Value = 1.001646801350442
$ws.Range("N22").Value = 1.001646801350442
$ws.Range("O22").Value = 1.001315009744215
$ws.Range("P22").Value = 0.9993787309659896
$ws.Range("Q22").Value = 0.9993787309659895
$ws.Range("R22").Value = 0.9982446957737635
$ws.Range("S22").Value = 0.9982446957737635
$ws.Range("T22").Value = 1.000367179500231

# Row 23: HexGrid-90degTilt15degRes
$ws.Range("C23").Value = 0.9999566947849777
$ws.Range("D23").Value = 1.014846105841701
$ws.Range("E23").Value = 1.002547176188888
$ws.Range("F23").Value = 0.97280487670134
$ws.Range("G23").Value = 0.97280487670134
$ws.Range("H23").Value = 1.003336163301229
$ws.Range("I23").Value = 1.003336163301229
$ws.Range("J23").Value = 1.011040282410761
$ws.Range("K23").Value = 0.97280487670134
$ws.Range("L23").Value = 1.011040282410761
$ws.Range("M23").Value = 1.007188222855995
$ws.Range("N23").Value = 1.007188222855995
$ws.Range("O23").Value = 1.005641207300292
$ws.Range("P23").Value = 0.9957271074711102
$ws.Range("Q23").Value = 0.9957271074711102
$ws.Range("R23").Value = 0.9899965497786676
$ws.Range("S23").Value = 0.9899965497786676
$ws.Range("T23").Value = 1.000755216538149

Write-Output $ws.UsedRange.Address()